$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("A8").Value = "Moodlight Quartet"
$ws.Range("A8").Font.Bold = $true

$ws.Range("B8").Value = "An RGB LED that reacts to light and motion in four different ways."

$ws.Range("B8").Select()
